# ランサーズ (sheet1) gets refreshed with a new scrape pulled at 2025-11-11 06:27:39.
# The freshest 2 postings are brand new, the rest re-appear (shifted up) from
# further down the previous day's list, and the tail of the old list (rows
# 10-20) is dropped entirely since the feed is capped at 8 listings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop any existing hyperlinks up front - row surgery below does not keep
# the <hyperlinks> collection in sync with the rows it references, so we
# rebuild it from scratch once the final F-column values are in place.
$ws.Hyperlinks.Delete()

# Remove the rows that fall off the end of the refreshed feed (old rows
# 10-20). This also shrinks the sheet's used range down to A1:H9.
$ws.Range("A10:A20").EntireRow.Delete()

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = "2025-11-11 06:27:39"
$ws.Range("B2").Value = "【急募】AIシステム構築!FirebaseとOpenAI活用の専門家募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5431299"
$ws.Range("G2").Value = 325
$ws.Range("H2").Value = "🔥AI,Ai"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = "2025-11-11 06:27:39"
$ws.Range("B3").Value = "【急募】大手保険会社向けスマホアプリ設計書作成依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5431609"
$ws.Range("G3").Value = 95
$ws.Range("H3").Value = "★スマホアプリ ◇アプリ"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = "2025-11-11 06:27:39"
$ws.Range("B4").Value = "【急募】知的財産関連システムの開発パートナー募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5431547"
$ws.Range("G4").Value = 90
$ws.Range("H4").Value = "◆開発"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = "2025-11-11 06:27:39"
$ws.Range("B5").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "~ 5,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = "◇管理"

# --- Row 6 (no skill digest in the refreshed feed) ------------------------
$ws.Range("A6").Value = "2025-11-11 06:27:39"
$ws.Range("B6").Value = "【急募】Laravel12でFortifyを使った2段階認証システムの制作"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5431508"
$ws.Range("G6").Value = 33
$ws.Range("H6").ClearContents()

# --- Row 7 (no skill digest) ----------------------------------------------
$ws.Range("A7").Value = "2025-11-11 06:27:39"
$ws.Range("B7").Value = "【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5431322"
$ws.Range("G7").Value = 25
$ws.Range("H7").ClearContents()

# --- Row 8 (no skill digest) ----------------------------------------------
$ws.Range("A8").Value = "2025-11-11 06:27:39"
$ws.Range("B8").Value = "AWS環境からAWS環境ヘの新規構築"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5431069"
$ws.Range("G8").Value = 18
$ws.Range("H8").ClearContents()

# --- Row 9 (no skill digest) ----------------------------------------------
$ws.Range("A9").Value = "2025-11-11 06:27:39"
$ws.Range("B9").Value = "EAの作成"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5431276"
$ws.Range("G9").Value = 10
$ws.Range("H9").ClearContents()

# Rebuild the hyperlinks for the URL column, F2:F9, in row order.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5431299")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5431609")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5431547")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5431508")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5431322")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5431069")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5431276")

# Column widths: B 52 -> 39, H 17 -> 14. ColumnWidth uses "characters of the
# default font" units which Excel pads by ~0.8333 vs the raw OOXML <col
# width>, so back that padding out to land on the exact target widths.
$ws.Columns.Item(2).ColumnWidth = 38.166666666666664
$ws.Columns.Item(8).ColumnWidth = 13.166666666666666

Write-Output "done"
